$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (shifts the existing Code/Description/Definition
# columns from A:C to B:D) and add the new "Version" column.
$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "Version"

# The version values ("1.0") must be stored as text, not as the number 1.
# Force text storage via a temporary text number format, then drop the
# format again so the cells keep their default (General) style.
$verRange = $ws.Range("A2:A6")
$verRange.NumberFormat = "@"
$verRange.Value = "1.0"
$verRange.ClearFormats()
